# Standard User name change in Contacts test data files - 6th Mar 2024
#
# Renames the Standard User contact name on the "Users" sheet from
# "Drew Koecher" to "Ayati Arvind", and updates the active sheet /
# selection state to reflect where the author left off editing
# (Users!D4 active, tabSelected; Contact!D12 selected, not the active tab).

$wb = $excel.ActiveWorkbook

# --- Data change -----------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Ayati Arvind"

# --- View / selection state -------------------------------------------
$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Range("D12").Select() | Out-Null

$wsUsers.Activate() | Out-Null
$wsUsers.Range("D4").Select() | Out-Null
